# Update "Förändrad" (Changed) date column C for rows 2-6 from 45183 to 45184
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 45184
$ws.Range("C3").Value = 45184
$ws.Range("C4").Value = 45184
$ws.Range("C5").Value = 45184
$ws.Range("C6").Value = 45184
